$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 21; everything currently at row 21 onward
# (including formatting) shifts down to make room.
$ws.Rows("21:23").Insert()

# Row 21: Camote, 1a (guarda)
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C21").Value = "Arica y Parinacota"
$ws.Range("D21").Value = "2023-05-25"
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 100112045
$ws.Range("G21").Value = "Zapallo"
$ws.Range("H21").Value = "Camote"
$ws.Range("I21").Value = "1a (guarda)"
$ws.Range("J21").Value = 700
$ws.Range("K21").Value = 370
$ws.Range("L21").Value = 390
$ws.Range("M21").Value = 380
$ws.Range("N21").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O21").Value = "Región de O'Higgins"
$ws.Range("P21").Value = 380
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"

# Row 22: Camote, 2a (guarda)
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = "2023-05-25"
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 100112045
$ws.Range("G22").Value = "Zapallo"
$ws.Range("H22").Value = "Camote"
$ws.Range("I22").Value = "2a (guarda)"
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 360
$ws.Range("M22").Value = 355
$ws.Range("N22").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O22").Value = "Región de O'Higgins"
$ws.Range("P22").Value = 355
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"

# Row 23: Camote, 3a (guarda)
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = "2023-05-25"
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112045
$ws.Range("G23").Value = "Zapallo"
$ws.Range("H23").Value = "Camote"
$ws.Range("I23").Value = "3a (guarda)"
$ws.Range("J23").Value = 800
$ws.Range("K23").Value = 330
$ws.Range("L23").Value = 340
$ws.Range("M23").Value = 335
$ws.Range("N23").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O23").Value = "Región de O'Higgins"
$ws.Range("P23").Value = 335
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
